$d = $word.ActiveDocument

# 1. Title/heading text changed in two places (same replacement both times)
$d.Content.Find.Execute("Play Betty Bonkers for Free: Read Our Review", $true, $false, $false, $false, $false, $true, 1, $false, "Play Betty Bonkers for free", 2)

# 2. Bullet list items under "What we like"
$d.Content.Find.Execute("Multiple bonus features and free spins", $true, $false, $false, $false, $false, $true, 1, $false, "Exciting gameplay features and bonuses", 2)
$d.Content.Find.Execute("Maximum win of over 10,000x the stake", $true, $false, $false, $false, $false, $true, 1, $false, "High maximum win potential", 2)
$d.Content.Find.Execute("Retro design with exceptional graphics and sound", $true, $false, $false, $false, $false, $true, 1, $false, "Retro design with exceptional graphics", 2)
$d.Content.Find.Execute("Wild and Scatter symbols for bigger payouts", $true, $false, $false, $false, $false, $true, 1, $false, "Buy option for quick access to bonuses", 2)

# 3. Bullet list items under "What we don't like"
$d.Content.Find.Execute("High volatility may not suit all players", $true, $false, $false, $false, $false, $true, 1, $false, "Limited number of paylines", 2)
$d.Content.Find.Execute("Only five fixed paylines may limit gameplay options", $true, $false, $false, $false, $false, $true, 1, $false, "High volatility may result in frequent losing streaks", 2)

# 4. Final italic summary paragraph
$d.Content.Find.Execute("Discover everything you need to know about Betty Bonkers online slot game, including its features, bonuses, RTP, and where to play it for free.", $true, $false, $false, $false, $false, $true, 1, $false, "Read our review of Betty Bonkers and play this exciting slot game for free.", 2)
